{"js": "// 1. Exercise 1: change the file password from 1234 to 5678.\n{\n  const results = context.document.body.search(\". The file password is 1234.\", { matchCase: true });\n  await context.sync();\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\". The file password is 5678.\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 2. Exercise 1: append a sentence about the MSF export password.\n{\n  const results = context.document.body.search(\"Export the filtered data with the anonymous export for MSF.\", { matchCase: true });\n  await context.sync();\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"Export the filtered data with the anonymous export for MSF. If you don\\u2019t change the export password, the password is 605637.\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// 3. Exercise 2 (Admission / hospitalisation): rework the sentence to reference\n//    the \"Linelist patients\" sheet and fix the spelling of \"hospitalisation\".\n{\n  const results = context.document.body.search(\"section, just before the\", { matchCase: true });\n  await context.sync();\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"section of\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  const hospResults = context.document.body.search(\"hospitalization\", { matchCase: true });\n  await context.sync();\n  hospResults.load(\"items\");\n  await context.sync();\n  if (hospResults.items.length > 0) {\n    const anchor = hospResults.items[0];\n    anchor.insertText(\"Linelist patients sheet, just before the \", Word.InsertLocation.before);\n    await context.sync();\n    anchor.insertText(\"hospitalisation\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  const newPhrase = context.document.body.search(\"Linelist patients\", { matchCase: true });\n  await context.sync();\n  newPhrase.load(\"items\");\n  await context.sync();\n  if (newPhrase.items.length > 0) {\n    newPhrase.items[0].font.italic = true;\n    await context.sync();\n  }\n}\n\n// 4. Exercise 2 (Vaccination): rework the sentence to reference the\n//    \"Linelist patients\" sheet as well.\n{\n  const results = context.document.body.search(\"section, just after\", { matchCase: true });\n  await context.sync();\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"section of\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  const vaccResults = context.document.body.search(\"Vaccination against measles\", { matchCase: true });\n  await context.sync();\n  vaccResults.load(\"items\");\n  await context.sync();\n  if (vaccResults.items.length > 0) {\n    vaccResults.items[0].insertText(\"Linelist patients sheet, just after \", Word.InsertLocation.before);\n    await context.sync();\n  }\n\n  const newPhrase = context.document.body.search(\"Linelist patients\", { matchCase: true });\n  await context.sync();\n  newPhrase.load(\"items\");\n  await context.sync();\n  if (newPhrase.items.length > 1) {\n    newPhrase.items[1].font.italic = true;\n    await context.sync();\n  }\n}\n\n// 5. Exercise 2: the first import_linelist_en.xlsb reference becomes .xlsx\n//    (the file used in Exercise 3 stays .xlsb).\n{\n  const results = context.document.body.search(\"import_linelist_en.xlsb\", { matchCase: true });\n  await context.sync();\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"import_linelist_en.xlsx\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 6. Exercise 3: update the temporal table instructions.\n{\n  const results = context.document.body.search(\n    \"In the analyses sheet, add a temporal table showing the evolution of discharge types by notification date. Add a percentage option in row, and do not add a total.\",\n    { matchCase: true }\n  );\n  await context.sync();\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"In the analyses sheet, add a temporal table showing the evolution of type of discharge by notification date. Add a percentage option in row, and add a total.\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Exercise 1: change the file password from 1234 to 5678.\n$r = $d.Content\n$r.Find.Execute(\". The file password is 1234.\", $false, $false, $false, $false, $false, $true, 1, $false, \". The file password is 5678.\", 2) | Out-Null\n\n# 2. Exercise 1: append a sentence about the MSF export password.\n$r = $d.Content\n$r.Find.Execute(\"Export the filtered data with the anonymous export for MSF.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Export the filtered data with the anonymous export for MSF. If you don\" + [char]0x2019 + \"t change the export password, the password is 605637.\", 2) | Out-Null\n\n# 3. Exercise 2 (Admission / hospitalisation): rework the sentence to reference\n#    the \"Linelist patients\" sheet and fix the spelling of \"hospitalisation\".\n$r = $d.Content\n$r.Find.Execute(\"section, just before the\", $false, $false, $false, $false, $false, $true, 1, $false, \"section of\", 2) | Out-Null\n\n$r = $d.Content\n$r.Find.Text = \"hospitalization\"\n$r.Find.Execute() | Out-Null\n$r.InsertBefore(\"Linelist patients sheet, just before the \")\n\n$r = $d.Content\n$r.Find.Execute(\"hospitalization\", $false, $false, $false, $false, $false, $true, 1, $false, \"hospitalisation\", 2) | Out-Null\n\n$r = $d.Content\n$r.Find.Text = \"Linelist patients\"\n$r.Find.Execute() | Out-Null\n$r.Italic = 1\n\n# 4. Exercise 2 (Vaccination): rework the sentence to reference the\n#    \"Linelist patients\" sheet as well.\n$r = $d.Content\n$r.Find.Execute(\"section, just after\", $false, $false, $false, $false, $false, $true, 1, $false, \"section of\", 2) | Out-Null\n\n$r = $d.Content\n$r.Find.Text = \"Vaccination against measles\"\n$r.Find.Execute() | Out-Null\n$r.InsertBefore(\"Linelist patients sheet, just after \")\n\n$r = $d.Content\n$r.Find.Text = \"Linelist patients\"\ndo {\n    $r.Find.Execute() | Out-Null\n} while ($r.Find.Found -and $r.Italic -eq 1)\nif ($r.Find.Found) {\n    $r.Italic = 1\n}\n\n# 5. Exercise 2: the first import_linelist_en.xlsb reference becomes .xlsx\n#    (the file used in Exercise 3 stays .xlsb).\n$r = $d.Content\n$r.Find.Text = \"import_linelist_en.xlsb\"\n$r.Find.Execute() | Out-Null\n$r.Text = \"import_linelist_en.xlsx\"\n\n# 6. Exercise 3: update the temporal table instructions.\n$r = $d.Content\n$r.Find.Execute(\"In the analyses sheet, add a temporal table showing the evolution of discharge types by notification date. Add a percentage option in row, and do not add a total.\", $false, $false, $false, $false, $false, $true, 1, $false, \"In the analyses sheet, add a temporal table showing the evolution of type of discharge by notification date. Add a percentage option in row, and add a total.\", 2) | Out-Null\n"}
